# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" sheet (copied from the "2022-Q2" sheet as a
# template, positioned immediately before it) with fresh Q3 figures, and
# records the new quarter in the "总计" (totals) summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "2022-Q3" worksheet by duplicating "2022-Q2" (same
#    layout/styling) and dropping it in right before "2022-Q2" so the
#    tab order becomes: 总计, 2022-Q3, 2022-Q2, 2021-Q4, 2021-Q3.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Fill in the Q3 figures. The percentage/number-like columns are
#    stored as text in this workbook, so a leading apostrophe forces
#    text entry, and resetting the style back to Normal afterwards
#    drops the quote-prefix formatting Excel would otherwise apply.
# ---------------------------------------------------------------------
$q3.Range("C2").Value = "华泰柏瑞亚洲领导企业混合（QDII）"

$q3.Range("D2").Value = "'0.36"
$q3.Range("D2").Style = "Normal"

$q3.Range("E2").Value = "'93.44"
$q3.Range("E2").Style = "Normal"

$q3.Range("F2").Value = "'5.06"
$q3.Range("F2").Style = "Normal"

$q3.Range("G2").Value = "'0.0182"
$q3.Range("G2").Style = "Normal"

$q3.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: insert a new row for 2021-Q3 at
#    the bottom and shift the quarter labels up so row 2 now reports
#    the brand-new 2022-Q3 quarter.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(5).Insert()
$total.Range("A4").Copy($total.Range("A5"))

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.02

$total.Range("B3").Value = "2022-Q2"
$total.Range("D3").Value = 0.03

$total.Range("B4").Value = "2021-Q4"
$total.Range("D4").Value = 0.06

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q3"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.07000000000000001

# ---------------------------------------------------------------------
# 4) Restore the originally-selected tab (2021-Q3), since adding /
#    copying sheets shifts Excel's active-sheet focus.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Select()
